$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update SELL PRICE (column E) values ---
$newSellPrices = @{
    2  = 7000
    3  = 9000
    4  = 5000
    5  = 6000
    6  = 7000
    7  = 31500
    8  = 33000
    9  = 30000
    10 = 25500
    11 = 30000
    12 = 96000
    13 = 105000
    14 = 96000
    15 = 111000
    16 = 105000
    17 = 325000
    18 = 350000
    19 = 275000
    20 = 275000
    21 = 350000
    22 = 300000
    23 = 1400000
    24 = 1000000
    25 = 1200000
    26 = 5000000
    27 = 4375000
    28 = 15000000
}

foreach ($rowNum in $newSellPrices.Keys) {
    $ws.Cells.Item($rowNum, 5).Value = $newSellPrices[$rowNum]
}

# --- Reset custom row heights back to the sheet default on the rows that ---
# --- previously wrapped text (rows 6, 7, 8, 21, 24) ---
$rowsToAutoFit = @(6, 7, 8, 21, 24)
foreach ($rowNum in $rowsToAutoFit) {
    $ws.Rows.Item($rowNum).AutoFit()
}

# --- Update the saved view state: scrolled-to column and active selection ---
$ws.Activate()
$ws.Range("S14").Select()
